$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the "data period" text shown in column B (rows 2-9) for the new
#    APS data release: "Apr 2024 - Mar 2025" -> "Jul 2024 - Jun 2025"
# ---------------------------------------------------------------------------
$rngB = $ws.Range("B2:B9")
$rngB.Value = "Jul 2024 - Jun 2025"

# ---------------------------------------------------------------------------
# 2. Update the destinations methodology text (column E, rows 20 & 21) so the
#    date ranges quoted move on by one academic year.
# ---------------------------------------------------------------------------
$e20 = "Destination measures show the percentage of students going to or remaining in an education, apprenticeship or employment destination in the academic year after completing Key Stage 4 studies (usually aged between 14 to 16). The cohort of learners used in the metrics here completed in AY22/23.`nA sustained destination is a count of young people recorded as having sustained participation (education and employment) for a 6 month period in the destination year.`nThis means attending for all of the first two terms of the academic year (e.g. October 2023 to March 2024) at one or more education providers; spending 5 of the 6 months in employment or a combination of the two.`nA sustained apprenticeship is recorded when 6 months continuous participation is recorded at any point in the destination year (between August 2023 and July 2024).`nNot recorded includes pupils who were captured in the destination source data but who failed to meet the sustained participation criteria.`nUnknown (activity not captured): The student was not found to have any participation in education, apprenticeship or employment nor recorded as receiving out-of-work benefits at any point in the year. This also includes not being recorded by their Local Authority as NEET (not engaged in education, employment or training)."

$e21 = "Destination measures show the percentage of students going to or remaining in an education, apprenticeship or employment destination in the academic year after completing Key Stage 5 studies (usually aged 18). The cohort of learners used in the metrics here completed in AY22/23.`nA sustained destination is a count of young people recorded as having sustained participation (education and employment) for a 6 month period in the destination year.`nThis means attending for all of the first two terms of the academic year (e.g. October 2023 to March 2024) at one or more education providers; spending 5 of the 6 months in employment or a combination of the two.`nA sustained apprenticeship is recorded when 6 months continuous participation is recorded at any point in the destination year (between August 2023 and July 2024).`nNot recorded includes pupils who were captured in the destination source data but who failed to meet the sustained participation criteria.`nUnknown (activity not captured): The student was not found to have any participation in education, apprenticeship or employment nor recorded as receiving out-of-work benefits at any point in the year. This also includes not being recorded by their Local Authority as NEET (not engaged in education, employment or training)."

$ws.Range("E20").Value = $e20
$ws.Range("E21").Value = $e21

# ---------------------------------------------------------------------------
# 3. Re-style the column-B data-period cells: smaller, muted Verdana font
#    instead of the previous Arial, and drop the explicit wrap/valign so the
#    cell falls back to the default alignment.
# ---------------------------------------------------------------------------
$rngB.ClearFormats()
$rngB.Font.Size = 7
$rngB.Font.Name = "Verdana"
$rngB.Font.Color = 5724003

# ---------------------------------------------------------------------------
# 4. Reset the view: scroll back to the top of the sheet and select E2
#    instead of the previous F20 selection / A20 scroll position.
# ---------------------------------------------------------------------------
$ws.Range("E2").Select()
